# Update "想去人数" (want-to-go count) values for two events that appear
# on both the "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 7691
    $ws.Range("F16").Value = 470
}
